# feat: add sequence field
#
# Replaces the "LVL" / "levelValue" parameter-level column with a new
# "#" / "sequence" column in the request- and response-parameter tables,
# and renames the stray "欄位名稱及說明" header cells (left over from an
# older layout) to the already-used "備註" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Parameter table column headers (row 11 "otherParameter" table) ---
# Column F header was "欄位名稱及說明" ("field name and description"); it
# should read "備註" ("remarks"), matching the other tables' wording.
$ws.Range("F11").Value = "備註"

# --- Request parameter table (rows 15-16) ---
# "LVL" -> "#" ; last column header "欄位名稱及說明" -> "備註".
$ws.Range("A15").Value = "#"
$ws.Range("F15").Value = "備註"
# Placeholder row: levelValue -> sequence, with the left-aligned text style
# (same numFmt "@" look as before, just horizontally left aligned) that the
# sibling response-parameter placeholder cell already uses.
$ws.Range("A16").Value = '${requestParameter.sequence}'
$ws.Range("A16").HorizontalAlignment = -4131

# --- Response parameter table (rows 22-23) ---
$ws.Range("A22").Value = "#"
$ws.Range("E22").Value = "備註"
$ws.Range("A23").Value = '${responseParameter.sequence}'

# --- Sheet view: scroll position / selection moved to the edited rows ---
$ws.Range("B8:F8").Select
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 7
